$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old rows beyond row 4 since the new data is more compact
$ws.Range("A5:A14").ClearContents()

# Update rows 2-4 with the new consolidated tuple-style strings
$ws.Range("A2").Value = "('Bat', ['Token Creature — Bat', 'Flying', '1/1'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '4/4'])"
$ws.Range("A4").Value = "('Elephant', ['Token Creature — Elephant', '3/3'])"
